$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Hours Report section - add a new bullet after "Modified the
# output of the table..." describing the Hours Generation change.
# ------------------------------------------------------------------
$pOutputTable = $d.Paragraphs.Item(5)
$pOutputTable.Range.InsertParagraphAfter()
$pNewHours = $d.Paragraphs.Item(6)
$pNewHours.Range.InsertAfter("Modified the Hours Generation to generate a start shift time and calculate the end shift time based on the hours worked (which is also randomly generated).  This is stored in the database as a datetime object (2025, 1, 1, 8, 45). (year, month, day, hour, minute).")
$pNewHours = $d.Paragraphs.Item(6)
$insertPt = $d.Range($pNewHours.Range.End - 1, $pNewHours.Range.End - 1)
$insertPt.InsertAfter("  Subsequently, the Report code had to be changed only to call on the correct column from the Database.")

# ------------------------------------------------------------------
# Step 2: Database Design section - remove the empty paragraph and the
# "ERD:" heading that used to sit right after the "Department Manager"
# bullet (it gets re-inserted further down, after a new bullet).
# ------------------------------------------------------------------
# After step 1, paragraph indices shifted by +1:
#   10 -> (empty)
#   11 -> "ERD:"
$pEmptyBeforeErd = $d.Paragraphs.Item(10)
$pEmptyBeforeErd.Range.Delete()
$pErdHeading = $d.Paragraphs.Item(10)
$pErdHeading.Range.Delete()

# ------------------------------------------------------------------
# Step 3: Replace the "Fixed a couple of the missing Foreign Key..."
# bullet text with the new "Refactored the hours table..." bullet.
# ------------------------------------------------------------------
$pFixedFk = $d.Paragraphs.Item(10)
$pFixedFk.Range.Text = "Refactored the hours table to have a StartShift/EndShift rather than a DateWorked to accurately reflect employees clocking in and out.  HoursWorked remains to more easily track the number of hours worked by the employee."

# ------------------------------------------------------------------
# Step 4: Re-insert (in order): blank paragraph, "ERD:" heading, and a
# new bullet carrying the original "Fixed a couple..." text.
# ------------------------------------------------------------------
$pRefactored = $d.Paragraphs.Item(10)
$pRefactored.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item(11)
$pBlank.Range.Style = "Normal"

$pBlank.Range.InsertParagraphAfter()
$pErdNew = $d.Paragraphs.Item(12)
$pErdNew.Range.Style = "Heading2"
$pErdNew.Range.InsertAfter("ERD:")

$pErdNew = $d.Paragraphs.Item(12)
$pErdNew.Range.InsertParagraphAfter()
$pFkNew = $d.Paragraphs.Item(13)
$pFkNew.Range.Style = "List Paragraph"
$numTemplate = $d.Paragraphs.Item(10).Range.ListFormat.ListTemplate
$pFkNew.Range.ListFormat.ApplyListTemplateWithLevel($numTemplate)
$pFkNew.Range.InsertAfter("Fixed a couple of the missing Foreign Key relationships on the ERD")

# ------------------------------------------------------------------
# Step 5: Add a final bullet after the "YearProduced" bullet describing
# the Hours Table change.
# ------------------------------------------------------------------
$pYearProduced = $d.Paragraphs.Item(14)
$pYearProduced.Range.InsertParagraphAfter()
$pHoursTable = $d.Paragraphs.Item(15)
$pHoursTable.Range.InsertAfter("Changed the Hours Table to match changes to the Database listed above.")

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Range.Text + "]")
}
